$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.913.99"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.549.84"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("E7").Value = "  -0.40%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.89%  "

$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "1.770.43"
$ws.Range("E12").Value = "  -0.11%  "

$ws.Range("D13").Value = "1.545.05"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("E14").Value = "  +0.90%  "

$ws.Range("E15").Value = "  +0.92%  "

$ws.Range("D16").Value = "26.911.37"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.64%  "

$ws.Range("E28").Value = "  +0.94%  "

$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.70%  "

$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.43%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.413.63"
$ws.Range("E34").Value = "  +3.50%  "

$ws.Range("E35").Value = "  +3.28%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E38").Value = "  +1.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.526"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.41%  "

$ws.Range("E43").Value = "  +3.11%  "

$ws.Range("E44").Value = "  +1.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").Value = "1.684.15"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.60%  "

$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
